$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset (old "RM 232" row 26, old "SC 92" row 28)
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Update individual data cells for rows 2-25 per the recorded changes
$ws.Range("E2").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("E5").Value = -5
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("D8").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("D12").Value = -14.1
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("F14").Value = 17.76
$ws.Range("F16").ClearContents()
$ws.Range("D17").Value = -14.7
$ws.Range("F17").ClearContents()
$ws.Range("D18").Value = -15.2
$ws.Range("D19").ClearContents()
$ws.Range("F19").Value = 17.81
$ws.Range("D20").ClearContents()
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = 16.81
$ws.Range("D23").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("F25").ClearContents()

# Update individual data cells for rows 26-33 (post row-shift) per the recorded changes
$ws.Range("F26").ClearContents()
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("E30").Value = -5.7
$ws.Range("F31").Value = 17.18
$ws.Range("C32").ClearContents()
